$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 8: section title "Moje poprawki"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Moje poprawki"

# ---------------------------------------------------------------------------
# Row 9: header row (duplicate of row 2) + new column J header
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "l. odpadów"
$ws.Range("B9").Font.Bold = $true

$ws.Range("C9").Value = "zdjęcia"
$ws.Range("D9").Value = "training"
$ws.Range("E9").Value = "test"

$ws.Range("F9").Value = "sum"
$ws.Range("F9").Font.Bold = $false
$ws.Range("F9").Font.Bold = $true
$ws.Range("F9").Font.Bold = $false

$ws.Range("G9").Value = "sum/zdjęcia ratio"
$ws.Range("G9").Font.Bold = $true

$ws.Range("H9").Value = "ile stopni obrotu"

$ws.Range("J9").Value = "wyszło:"
$ws.Range("J9").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# Row 10: PET (duplicate of row 3)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "PET"

$ws.Range("C10").Value = 2200
$ws.Range("C10").HorizontalAlignment = -4108

$ws.Range("D10").Value = 32400
$ws.Range("D10").HorizontalAlignment = -4108

$ws.Range("E10").Value = 3600
$ws.Range("E10").HorizontalAlignment = -4108

$ws.Range("B10").Formula = "=C10/40"
$ws.Range("B10").Font.Bold = $true
$ws.Range("B10").HorizontalAlignment = -4108

$ws.Range("F10").Formula = "=SUM(D10,E10)"
$ws.Range("F10").Font.Bold = $false
$ws.Range("F10").Font.Bold = $true
$ws.Range("F10").Font.Bold = $false
$ws.Range("F10").HorizontalAlignment = -4108

$ws.Range("G10").Formula = "=F10/C10"
$ws.Range("G10").Font.Bold = $true
$ws.Range("G10").HorizontalAlignment = -4108
$ws.Range("G10").NumberFormat = "0"

$ws.Range("H10").Formula = "=360/G10"
$ws.Range("H10").HorizontalAlignment = -4108
$ws.Range("H10").NumberFormat = "0"

$ws.Range("J10").Value = 35200

# ---------------------------------------------------------------------------
# Row 11: HDPE (duplicate of row 4)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "HDPE"

$ws.Range("C11").Value = 600
$ws.Range("C11").HorizontalAlignment = -4108

$ws.Range("D11").Value = 29700
$ws.Range("D11").HorizontalAlignment = -4108

$ws.Range("E11").Value = 3300
$ws.Range("E11").HorizontalAlignment = -4108

$ws.Range("B11").Formula = "=C11/40"
$ws.Range("B11").Font.Bold = $true
$ws.Range("B11").HorizontalAlignment = -4108

$ws.Range("F11").Formula = "=SUM(D11,E11)"
$ws.Range("F11").Font.Bold = $false
$ws.Range("F11").Font.Bold = $true
$ws.Range("F11").Font.Bold = $false
$ws.Range("F11").HorizontalAlignment = -4108

$ws.Range("G11").Formula = "=F11/C11"
$ws.Range("G11").Font.Bold = $true
$ws.Range("G11").HorizontalAlignment = -4108

$ws.Range("H11").Formula = "=360/G11"
$ws.Range("H11").NumberFormat = "0"
$ws.Range("H11").HorizontalAlignment = -4108

$ws.Range("J11").Value = 33000

# ---------------------------------------------------------------------------
# Row 12: PP (duplicate of row 5)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "PP"

$ws.Range("C12").Value = 640
$ws.Range("C12").HorizontalAlignment = -4108

$ws.Range("D12").Value = 29952
$ws.Range("D12").HorizontalAlignment = -4108

$ws.Range("E12").Value = 3328
$ws.Range("E12").HorizontalAlignment = -4108

$ws.Range("B12").Formula = "=C12/40"
$ws.Range("B12").Font.Bold = $true
$ws.Range("B12").HorizontalAlignment = -4108

$ws.Range("F12").Formula = "=SUM(D12,E12)"
$ws.Range("F12").Font.Bold = $false
$ws.Range("F12").Font.Bold = $true
$ws.Range("F12").Font.Bold = $false
$ws.Range("F12").HorizontalAlignment = -4108

$ws.Range("G12").Formula = "=F12/C12"
$ws.Range("G12").Font.Bold = $true
$ws.Range("G12").HorizontalAlignment = -4108
$ws.Range("G12").NumberFormat = "0"

$ws.Range("H12").Formula = "=360/G12"
$ws.Range("H12").NumberFormat = "0"
$ws.Range("H12").HorizontalAlignment = -4108

$ws.Range("J12").Value = 33280
$ws.Range("J12").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# Row 13: PS (duplicate of row 6)
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "PS"

$ws.Range("C13").Value = 520
$ws.Range("C13").HorizontalAlignment = -4108

$ws.Range("D13").Value = 28080
$ws.Range("D13").HorizontalAlignment = -4108

$ws.Range("E13").Formula = "=D13*0.1"
$ws.Range("E13").HorizontalAlignment = -4108

$ws.Range("B13").Formula = "=C13/40"
$ws.Range("B13").Font.Bold = $true
$ws.Range("B13").HorizontalAlignment = -4108

$ws.Range("F13").Formula = "=SUM(D13,E13)"
$ws.Range("F13").Font.Bold = $false
$ws.Range("F13").Font.Bold = $true
$ws.Range("F13").Font.Bold = $false
$ws.Range("F13").HorizontalAlignment = -4108

$ws.Range("G13").Formula = "=F13/C13"
$ws.Range("G13").Font.Bold = $true
$ws.Range("G13").HorizontalAlignment = -4108
$ws.Range("G13").NumberFormat = "0.0"

$ws.Range("H13").Formula = "=360/G13"
$ws.Range("H13").NumberFormat = "0"
$ws.Range("H13").HorizontalAlignment = -4108

$ws.Range("J13").Value = 31200

# ---------------------------------------------------------------------------
# Row 15: lone formatted (empty) cell
# ---------------------------------------------------------------------------
$ws.Range("C15").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Existing rows 4-6: G column restyle (s=3 -> s=4) + un-share G6
# ---------------------------------------------------------------------------
$ws.Range("G4").NumberFormat = "0"
$ws.Range("G5").NumberFormat = "0"
$ws.Range("G6").Formula = "=F6/C6"
$ws.Range("G6").NumberFormat = "0"

# ---------------------------------------------------------------------------
# View state: active cell selection
# ---------------------------------------------------------------------------
$ws.Range("I5").Select()
